# Applies the "Assignment2_Groups" edit:
#  - Splits several name runs so that each individual first/last name that
#    Word's proofer would flag is wrapped in <w:proofErr w:type="spellStart"/>
#    ... <w:proofErr w:type="spellEnd"/> markers (matching a real proofing pass).
#  - Folds "Abhi Nileshkumar Patel" into the "Xin Zhao, Neel Ajay Mahimkar"
#    bullet (as plain, non-bold text) and removes its old standalone bullet.
#
# Strategy: for each target paragraph, replace the *entire* paragraph range
# (Start..End, i.e. including the trailing paragraph mark) with a pkg:package
# InsertXML payload containing a single <w:p> whose runs are the desired
# content. When the target Range spans a full paragraph including its mark,
# InsertXML preserves the paragraph's own <w:pPr> (style/numbering/etc.) and
# simply swaps in the new run content - exactly what's needed here.

$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $d.Range($p.Range.Start, $p.Range.End)
    $xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" + `
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" + `
           "<pkg:xmlData>" + `
           "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
           "<w:body><w:p>" + $innerXml + "</w:p></w:body>" + `
           "</w:document>" + `
           "</pkg:xmlData></pkg:part></pkg:package>"
    $full.InsertXML($xml)
}

# --- Paragraph 3: "Bhavjot Pal, Samay Sehgal, Kannav Sethi" ---
$inner3 = `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Bhavjot</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> Pal, </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Samay</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> Sehgal, </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Kannav</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> Sethi</w:t></w:r>"
Replace-ParagraphXml 3 $inner3

# --- Paragraph 4: "Yiyuan Dong" ---
$inner4 = `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Yiyuan</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> Dong</w:t></w:r>"
Replace-ParagraphXml 4 $inner4

# --- Paragraph 5: "Huu Minh Phong Nguyen" (fr-FR on every run) ---
$frRpr = "<w:rPr><w:lang w:val=""fr-FR""/></w:rPr>"
$inner5 = `
  "<w:r>" + $frRpr + "<w:t xml:space=""preserve"">Huu Minh </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r>" + $frRpr + "<w:t>Phong</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r>" + $frRpr + "<w:t xml:space=""preserve""> Nguyen</w:t></w:r>"
Replace-ParagraphXml 5 $inner5

# --- Paragraph 6: "Xin Zhao, Neel Ajay Mahimkar" -> also appends
#     ", Abhi Nileshkumar Patel" (merged in from the old bullet 8) ---
$inner6 = `
  "<w:r><w:t xml:space=""preserve"">Xin Zhao, </w:t></w:r>" + `
  "<w:r><w:t xml:space=""preserve"">Neel Ajay </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Mahimkar</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve"">, </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Abhi</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Nileshkumar</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> Patel</w:t></w:r>"
Replace-ParagraphXml 6 $inner6

# --- Paragraph 7: "Dennis Audu, Orang Tang Enow, Mohammadsadegh Firouzi" ---
$inner7 = `
  "<w:r><w:t xml:space=""preserve"">Dennis </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Audu</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve"">, </w:t></w:r>" + `
  "<w:r><w:t xml:space=""preserve"">Orang Tang Enow, </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Mohammadsadegh</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>" + `
  "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" + `
  "<w:proofErr w:type=""spellStart""/>" + `
  "<w:r><w:t>Firouzi</w:t></w:r>" + `
  "<w:proofErr w:type=""spellEnd""/>"
Replace-ParagraphXml 7 $inner7

# --- Paragraph 8: old standalone "Abhi Nileshkumar Patel" bullet - remove it
#     entirely now that its text lives in paragraph 6. ---
$p8 = $d.Paragraphs(8)
$d.Range($p8.Range.Start, $p8.Range.End).Delete()
